$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Row 25: new entry - Rédaction, UML/rapport work
$ws.Range("A25").Value = (Get-Date -Year 2023 -Month 4 -Day 24 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B25").Value = "Rédaction"
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = "UML, rapport"

# Row 26: new entry - Implémentation, backend restructuring work
$ws.Range("A26").Value = (Get-Date -Year 2023 -Month 4 -Day 24 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B26").Value = "Implémentation"
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = "Installation de l'ORM Sequelize, refactor service - repository, amélioration système de modules"

$wb.Save()
